$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("H62").Value = 5562.9165
$ws.Range("I62").Value = 4930.8335
$ws.Range("J62").Value = 6195
$ws.Range("K62").Value = 4930.8335
$ws.Range("L62").Value = 6195
$ws.Range("M62").Value = -4306.8335
$ws.Range("N62").Value = -7443

# Row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("H65").Value = 5562.9165
$ws.Range("I65").Value = 4930.8335
$ws.Range("J65").Value = 6195
$ws.Range("K65").Value = 24654.1675
$ws.Range("L65").Value = 30975
$ws.Range("M65").Value = -21534.1675
$ws.Range("N65").Value = -37215

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 1507.6666
$ws.Range("I98").Value = 1531.8889
$ws.Range("J98").Value = 1435
$ws.Range("K98").Value = 1531.8889
$ws.Range("L98").Value = 1435
$ws.Range("M98").Value = -33.88889999999992
$ws.Range("N98").Value = -4431

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 1507.6666
$ws.Range("I122").Value = 1531.8889
$ws.Range("J122").Value = 1435
$ws.Range("K122").Value = 4595.6667
$ws.Range("L122").Value = 4305
$ws.Range("M122").Value = -2145.6667
$ws.Range("N122").Value = -9205

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = ""
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = 0

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2632.889
$ws.Range("I138").Value = 1898.6666
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 5695.9998
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -555.9997999999996
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 2500
$ws.Range("I2").Value = 2500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -2387

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 2105.6897
$ws.Range("I32").Value = 1823.75
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 1823.75
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -1536.75
$ws.Range("N32").Value = -10574

# Row 53: Metal Fatigue | Mythril Vambraces
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = ""

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 322.47058
$ws.Range("I97").Value = 311.375
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 311.375
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = 184.625
$ws.Range("N97").Value = -1492

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 2636
$ws.Range("I110").Value = 2707.6667
$ws.Range("J110").Value = 2609.125
$ws.Range("K110").Value = 2707.6667
$ws.Range("L110").Value = 2609.125
$ws.Range("M110").Value = -662.6667000000002
$ws.Range("N110").Value = -6699.125

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 2500
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = -206

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 5170.6665
$ws.Range("I132").Value = 5170.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15511.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12981.9995

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 2500
$ws.Range("I3").Value = 2500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = -2386

# Row 15: Anutha Spatha | Bronze Spatha
$ws.Range("H15").Value = 29985.334
$ws.Range("I15").Value = 18000
$ws.Range("J15").Value = 35978
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 35978
$ws.Range("M15").Value = -17773
$ws.Range("N15").Value = -36432

# Row 19: Twice as Slice | Spiked Bronze Labrys
$ws.Range("H19").Value = 25980
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 25980
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 25980
$ws.Range("N19").Value = -26326

# Row 35: Lancers' Creed | Crowsbeak Hammer
$ws.Range("H35").Value = 15000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2397.923
$ws.Range("I105").Value = 2406.0833
$ws.Range("J105").Value = 2300
$ws.Range("K105").Value = 2406.0833
$ws.Range("L105").Value = 2300
$ws.Range("M105").Value = -659.0832999999998
$ws.Range("N105").Value = -5794

# Row 123: Archon Denied | High Durium Saw
$ws.Range("H123").Value = 126000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 126000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 126000
$ws.Range("N123").Value = -135800

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 8577.111000000001
$ws.Range("I134").Value = 8399.25
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 25197.75
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -22662.75
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 4674.5
$ws.Range("I58").Value = 3474.25
$ws.Range("J58").Value = 5874.75
$ws.Range("K58").Value = 3474.25
$ws.Range("L58").Value = 5874.75
$ws.Range("M58").Value = -3271.25
$ws.Range("N58").Value = -6280.75

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 3648.9
$ws.Range("I99").Value = 3387.6667
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 3387.6667
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -1889.6667
$ws.Range("N99").Value = -8996

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 2164.8
$ws.Range("I122").Value = 2164.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6494.400000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -4044.400000000001

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 3648.9
$ws.Range("I126").Value = 3387.6667
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 10163.0001
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -7693.000100000001
$ws.Range("N126").Value = -22940

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = ""
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = 0

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = ""
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = 0

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 4674.5
$ws.Range("I136").Value = 3474.25
$ws.Range("J136").Value = 5874.75
$ws.Range("K136").Value = 10422.75
$ws.Range("L136").Value = 17624.25
$ws.Range("M136").Value = -7872.75
$ws.Range("N136").Value = -22724.25

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry | Kukuru Powder
$ws.Range("H14").Value = 20262.2
$ws.Range("I14").Value = 20262.2
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 60786.60000000001
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -60613.60000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 937.625
$ws.Range("I113").Value = 937.625
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 937.625
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1232.375

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 10100
$ws.Range("I126").Value = 9533.333000000001
$ws.Range("J126").Value = 11800
$ws.Range("K126").Value = 28599.999
$ws.Range("L126").Value = 35400
$ws.Range("M126").Value = -26129.999
$ws.Range("N126").Value = -40340

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3512
$ws.Range("I132").Value = 3512
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10536
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8006

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 4199.4
$ws.Range("I16").Value = 4265.6665
$ws.Range("J16").Value = 4100
$ws.Range("K16").Value = 4265.6665
$ws.Range("L16").Value = 4100
$ws.Range("M16").Value = -4095.6665
$ws.Range("N16").Value = -4440

# Row 24: On Their Feet Again | Hard Leather Espadrilles
$ws.Range("H24").Value = 36000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 36000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 36000
$ws.Range("N24").Value = -36686

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 1795
$ws.Range("I81").Value = 1790
$ws.Range("J81").Value = 1800
$ws.Range("K81").Value = 3580
$ws.Range("L81").Value = 3600
$ws.Range("M81").Value = -2519
$ws.Range("N81").Value = -5722

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 1795
$ws.Range("I84").Value = 1790
$ws.Range("J84").Value = 1800
$ws.Range("K84").Value = 17900
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = -12596
$ws.Range("N84").Value = -28608

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2124.6155
$ws.Range("I122").Value = 2212
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 6636
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -4186
$ws.Range("N122").Value = -10400.0002

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1994
$ws.Range("I132").Value = 1994
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5982
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -3452
